$d = $word.ActiveDocument

# Fix the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-12-12 Thursday"

# Update each math expression cell in the table (20 rows x 5 columns)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "32+36="
$t.Cell(1, 2).Range.Text = "50-21="
$t.Cell(1, 3).Range.Text = "13+37="
$t.Cell(1, 4).Range.Text = "60+19="
$t.Cell(1, 5).Range.Text = "98-70="

$t.Cell(2, 1).Range.Text = "97-50="
$t.Cell(2, 2).Range.Text = "86-74="
$t.Cell(2, 3).Range.Text = "11+26="
$t.Cell(2, 4).Range.Text = "5+21="
$t.Cell(2, 5).Range.Text = "67+24="

$t.Cell(3, 1).Range.Text = "86-73="
$t.Cell(3, 2).Range.Text = "87-76="
$t.Cell(3, 3).Range.Text = "99-23="
$t.Cell(3, 4).Range.Text = "37+15="
$t.Cell(3, 5).Range.Text = "45-39="

$t.Cell(4, 1).Range.Text = "2-1="
$t.Cell(4, 2).Range.Text = "69-59="
$t.Cell(4, 3).Range.Text = "80-13="
$t.Cell(4, 4).Range.Text = "23+36="
$t.Cell(4, 5).Range.Text = "82-17="

$t.Cell(5, 1).Range.Text = "4+38="
$t.Cell(5, 2).Range.Text = "43+12="
$t.Cell(5, 3).Range.Text = "53-33="
$t.Cell(5, 4).Range.Text = "4+4="
$t.Cell(5, 5).Range.Text = "19-3="

$t.Cell(6, 1).Range.Text = "41+5="
$t.Cell(6, 2).Range.Text = "94-89="
$t.Cell(6, 3).Range.Text = "92-32="
$t.Cell(6, 4).Range.Text = "88-19="
$t.Cell(6, 5).Range.Text = "73-43="

$t.Cell(7, 1).Range.Text = "96-31="
$t.Cell(7, 2).Range.Text = "85-71="
$t.Cell(7, 3).Range.Text = "90-44="
$t.Cell(7, 4).Range.Text = "80-1="
$t.Cell(7, 5).Range.Text = "7+32="

$t.Cell(8, 1).Range.Text = "57+41="
$t.Cell(8, 2).Range.Text = "7+86="
$t.Cell(8, 3).Range.Text = "21+65="
$t.Cell(8, 4).Range.Text = "67-47="
$t.Cell(8, 5).Range.Text = "30+63="

$t.Cell(9, 1).Range.Text = "79-39="
$t.Cell(9, 2).Range.Text = "50+48="
$t.Cell(9, 3).Range.Text = "56-20="
$t.Cell(9, 4).Range.Text = "66+14="
$t.Cell(9, 5).Range.Text = "50+11="

$t.Cell(10, 1).Range.Text = "58-31="
$t.Cell(10, 2).Range.Text = "15+54="
$t.Cell(10, 3).Range.Text = "96-18="
$t.Cell(10, 4).Range.Text = "73-56="
$t.Cell(10, 5).Range.Text = "36+16="

$t.Cell(11, 1).Range.Text = "80-74="
$t.Cell(11, 2).Range.Text = "84+14="
$t.Cell(11, 3).Range.Text = "15+33="
$t.Cell(11, 4).Range.Text = "85-82="
$t.Cell(11, 5).Range.Text = "84-48="

$t.Cell(12, 1).Range.Text = "31+19="
$t.Cell(12, 2).Range.Text = "18+58="
$t.Cell(12, 3).Range.Text = "88+8="
$t.Cell(12, 4).Range.Text = "7-1="
$t.Cell(12, 5).Range.Text = "5+10="

$t.Cell(13, 1).Range.Text = "20+25="
$t.Cell(13, 2).Range.Text = "35+50="
$t.Cell(13, 3).Range.Text = "50-27="
$t.Cell(13, 4).Range.Text = "30+39="
$t.Cell(13, 5).Range.Text = "51-15="

$t.Cell(14, 1).Range.Text = "19+43="
$t.Cell(14, 2).Range.Text = "50+47="
$t.Cell(14, 3).Range.Text = "98-12="
$t.Cell(14, 4).Range.Text = "8+59="
$t.Cell(14, 5).Range.Text = "67-16="

$t.Cell(15, 1).Range.Text = "41+37="
$t.Cell(15, 2).Range.Text = "38-25="
$t.Cell(15, 3).Range.Text = "65+31="
$t.Cell(15, 4).Range.Text = "17+1="
$t.Cell(15, 5).Range.Text = "8+1="

$t.Cell(16, 1).Range.Text = "35-14="
$t.Cell(16, 2).Range.Text = "23+39="
$t.Cell(16, 3).Range.Text = "44-27="
$t.Cell(16, 4).Range.Text = "17+42="
$t.Cell(16, 5).Range.Text = "86-81="

$t.Cell(17, 1).Range.Text = "14+30="
$t.Cell(17, 2).Range.Text = "37+33="
$t.Cell(17, 3).Range.Text = "34-18="
$t.Cell(17, 4).Range.Text = "57+16="
$t.Cell(17, 5).Range.Text = "24+15="

$t.Cell(18, 1).Range.Text = "34-23="
$t.Cell(18, 2).Range.Text = "78-33="
$t.Cell(18, 3).Range.Text = "43+19="
$t.Cell(18, 4).Range.Text = "8+10="
$t.Cell(18, 5).Range.Text = "39-27="

$t.Cell(19, 1).Range.Text = "20-16="
$t.Cell(19, 2).Range.Text = "81-22="
$t.Cell(19, 3).Range.Text = "10+20="
$t.Cell(19, 4).Range.Text = "98-18="
$t.Cell(19, 5).Range.Text = "90-47="

$t.Cell(20, 1).Range.Text = "6+54="
$t.Cell(20, 2).Range.Text = "92-82="
$t.Cell(20, 3).Range.Text = "83-8="
$t.Cell(20, 4).Range.Text = "22+8="
$t.Cell(20, 5).Range.Text = "21-11="
